$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row continuing the existing Day_Number / Date series
# (row 41 held day 40 / 2019-02-09, so row 42 is day 41 / 2019-02-10).
$ws.Range("A42").Value = 41
$ws.Range("B42").Value = 43506

# Copy the formatting (number format, borders, font, alignment) from the
# row above so the new row matches the existing style pattern exactly.
$ws.Range("A41:B41").Copy()
$ws.Range("A42:B42").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Match the row height used by the rest of the data rows.
$ws.Rows.Item(42).RowHeight = $ws.Rows.Item(41).RowHeight

# Match the author's recorded selection after the edit.
$ws.Range("D39").Select()
